# Add a new "2022" column (S) to the manufacturing GVA table, mirroring the
# existing year columns (D..R = 2007..2021), and revise a few of the most
# recent years' figures (P, Q, R) for rows 4 and 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column S formatting, cloned from the neighbouring column R ---
# Row 2 is a thin spacer row above the header; row 3 is the year header row.
$ws.Range("R2").Copy($ws.Range("S2"))
$ws.Range("R3").Copy($ws.Range("S3"))
$ws.Range("R4").Copy($ws.Range("S4"))
$ws.Range("R5").Copy($ws.Range("S5"))

$excel.CutCopyMode = $false

# --- New column S values ---
$ws.Range("S3").Value = 2022
$ws.Range("S4").Value = 13.6
$ws.Range("S5").Value = 20

# --- Revised figures for 2019-2021 (columns P, Q, R) ---
$ws.Range("P4").Value = 13.7
$ws.Range("Q4").Value = 13.1
$ws.Range("R4").Value = 11.8

$ws.Range("P5").Value = 13.6
$ws.Range("Q5").Value = 12.5
$ws.Range("R5").Value = 13.5

# --- Selection now sits on the new column's spacer cell ---
$ws.Range("S2").Select()
